{"js": "// The document is a date heading followed by a 5-column practice table of\n// \"two-digit \u00f7 one-digit\" problems. The edit swaps the date and every\n// populated problem cell for a new value while leaving the paragraph /\n// table structure (fonts, alignment, row/column counts) untouched.\n//\n// Because every non-blank paragraph in the document gets a new value, and\n// the number (and order) of non-blank paragraphs is unchanged, the safest\n// way to reproduce the edit is to walk the paragraphs in document order and\n// replace the text of each non-empty one with the corresponding new value,\n// using `insertText(..., Word.InsertLocation.replace)` on the paragraph so\n// existing run/paragraph formatting (fonts, size, justification) is kept.\n\nconst newValues = [\n  \"2025-10-20 Monday\",\n  \"98\u00f75=\",\n  \"32\u00f77=\",\n  \"98\u00f73=\",\n  \"59\u00f72=\",\n  \"12\u00f76=\",\n  \"90\u00f79=\",\n  \"52\u00f73=\",\n  \"23\u00f78=\",\n  \"39\u00f76=\",\n  \"62\u00f74=\",\n  \"24\u00f73=\",\n  \"13\u00f79=\",\n  \"84\u00f75=\",\n  \"82\u00f73=\",\n  \"36\u00f77=\",\n  \"20\u00f74=\",\n  \"13\u00f74=\",\n  \"61\u00f73=\",\n  \"33\u00f73=\",\n  \"26\u00f79=\",\n  \"35\u00f76=\",\n  \"11\u00f74=\",\n  \"80\u00f78=\",\n  \"72\u00f74=\",\n  \"36\u00f78=\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet valueIndex = 0;\nfor (let i = 0; i < paragraphs.items.length && valueIndex < newValues.length; i++) {\n  const paragraph = paragraphs.items[i];\n  if (paragraph.text !== \"\") {\n    paragraph.insertText(newValues[valueIndex], Word.InsertLocation.replace);\n    valueIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document is a date heading followed by a single 5-column, 20-row\n# practice table of \"two-digit \u00f7 one-digit\" problems (only every other row\n# holds visible text; the rest are spacer rows). The edit swaps the date\n# and every populated problem cell for a new value while leaving the\n# paragraph / table structure (fonts, alignment, row/column counts)\n# untouched. Assigning to `Range.Text` on a paragraph/cell range preserves\n# that range's existing run/paragraph formatting (fonts, size,\n# justification) while only swapping the visible text, which matches the\n# target edit exactly.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading (first paragraph in the document body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-10-20 Monday\"\n\n# 2) Table problem cells. Only rows 1, 5, 9, 13, 17 (1-based) contain text;\n#    the rows in between are blank spacer rows. Values below are listed in\n#    row-major reading order, matching the document order of the cells.\n$tbl = $d.Tables.Item(1)\n$contentRows = @(1, 5, 9, 13, 17)\n$newValues = @(\n    \"98\u00f75=\", \"32\u00f77=\", \"98\u00f73=\", \"59\u00f72=\", \"12\u00f76=\",\n    \"90\u00f79=\", \"52\u00f73=\", \"23\u00f78=\", \"39\u00f76=\", \"62\u00f74=\",\n    \"24\u00f73=\", \"13\u00f79=\", \"84\u00f75=\", \"82\u00f73=\", \"36\u00f77=\",\n    \"20\u00f74=\", \"13\u00f74=\", \"61\u00f73=\", \"33\u00f73=\", \"26\u00f79=\",\n    \"35\u00f76=\", \"11\u00f74=\", \"80\u00f78=\", \"72\u00f74=\", \"36\u00f78=\"\n)\n\n$i = 0\nforeach ($r in $contentRows) {\n    for ($c = 1; $c -le 5; $c++) {\n        $tbl.Cell($r, $c).Range.Text = $newValues[$i]\n        $i = $i + 1\n    }\n}\n"}
